# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.691.12"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.471.19"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.552"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.45%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0892"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "2.853.64"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.07%  "
$ws.Range("D16").Value = "2.469.93"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("E17").Value = "  +3.52%  "
$ws.Range("D18").Value = "41.655.68"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "0.0₃0977"
$ws.Range("E19").Value = "  +5.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "1.968.30"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.48%  "
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("D48").Value = "2.706.59"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.27%  "
